$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values to D1 and E1 (extends used range to A1:E6)
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3

# Update the selected cell/range to match the author's final selection
$ws.Range("F3").Select()
